# Downgrade Python version to 3.11
# The source sheet listed per-image metadata (Filename/Format/Mode/Size/
# Width/Height/Extracted Text) for a single image in columns A:G.
# Collapse it down to a single "Extracted Text" column, with the OCR text
# that used to be crammed (with embedded newlines) into G2 split out one
# line per row below the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the Filename/Format/Mode/Size/Width/Height columns entirely (shifts
# nothing since they're the trailing columns next to A) so only column A
# remains.
$ws.Range("B:G").Delete()

# Header
$ws.Range("A1").Value = "Extracted Text"

# Former G2 text, split into one row per line.
$ws.Range("A2").Value = "ss 220854 30 pom poe vawn/xoanten"
$ws.Range("A3").Value = "oo? «SY8stUzESBOSE “took xan REA Or swntO ant Soepore SEE T=! ta feud 4a od o00t"
$ws.Range("A4").Value = "ETAT >055 30 tng T/G0S 295 tdang S05 ag"
$ws.Range("A5").Value = "Teter Sc/e0/o aa"
